$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 80
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 25
$ws.Range("B8").Value = 26000000

$ws.Range("C2").Interior.Color = 32768
